$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $text)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextCell "D2" "26.765.90"
Set-TextCell "D3" "1.566.71"
Set-TextCell "E3" "  +0.02%  "
Set-TextCell "D5" "206.40"
Set-TextCell "E5" "  -0.99%  "
Set-TextCell "E6" "  -2.43%  "
Set-TextCell "D8" "21.92"
Set-TextCell "E8" "  -0.81%  "
Set-TextCell "E9" "  -0.66%  "
Set-TextCell "E10" "  -1.36%  "
Set-TextCell "E11" "  -0.23%  "
Set-TextCell "D12" "1.789.28"
Set-TextCell "E12" "  +0.06%  "
Set-TextCell "D13" "1.565.35"
Set-TextCell "E13" "  -0.42%  "
Set-TextCell "E14" "  -2.43%  "
Set-TextCell "E15" "  -0.53%  "
Set-TextCell "D16" "61.52"
Set-TextCell "E16" "  -3.35%  "
Set-TextCell "D17" "26.788.90"
Set-TextCell "D18" "215.06"
Set-TextCell "E18" "  +0.93%  "
Set-TextCell "D19" "7.38"
Set-TextCell "D20" "0.0₃0678"
Set-TextCell "E20" "  -1.71%  "
Set-TextCell "E21" "  +0.16%  "
Set-TextCell "D22" "4.10"
Set-TextCell "E22" "  -0.10%  "
Set-TextCell "E23" "  -2.47%  "
Set-TextCell "E24" "  -1.60%  "
Set-TextCell "D25" "152.44"
Set-TextCell "E25" "  -0.58%  "
Set-TextCell "E26" "  +0.87%  "
Set-TextCell "D27" "14.94"
Set-TextCell "E27" "  -0.23%  "
Set-TextCell "E28" "  -0.03%  "
Set-TextCell "E29" "  -1.34%  "
Set-TextCell "B30" "PancakeSwap"
Set-TextCell "C30" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D30" "1.11"
Set-TextCell "E30" "  -3.41%  "
Set-TextCell "B31" "Hedera"
Set-TextCell "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D31" "0.0464"
Set-TextCell "E31" "  -1.18%  "
Set-TextCell "D32" "3.16"
Set-TextCell "E32" "  -1.32%  "
Set-TextCell "D33" "1.389.58"
Set-TextCell "E33" "  +1.00%  "
Set-TextCell "E34" "  -1.23%  "
Set-TextCell "E35" "  +0.60%  "
Set-TextCell "D37" "0.930"
Set-TextCell "E37" "  -2.98%  "
Set-TextCell "D39" "0.526"
Set-TextCell "E39" "  -0.83%  "
Set-TextCell "D40" "0.818"
Set-TextCell "E40" "  -0.57%  "
Set-TextCell "E41" "  +0.08%  "
Set-TextCell "D42" "0.989"
Set-TextCell "E42" "  +1.39%  "
Set-TextCell "D43" "1.80"
Set-TextCell "E43" "  -0.26%  "
Set-TextCell "E44" "  +1.73%  "
Set-TextCell "E45" "  +0.84%  "
Set-TextCell "D46" "63.35"
Set-TextCell "E46" "  -1.09%  "
Set-TextCell "D47" "1.702.21"
Set-TextCell "E47" "  +0.27%  "
Set-TextCell "D48" "85.63"
Set-TextCell "E48" "  +0.16%  "
Set-TextCell "D49" "0.0₇0990"
Set-TextCell "E49" "  -0.96%  "
Set-TextCell "D50" "0.0951"
Set-TextCell "E50" "  -0.73%  "
Set-TextCell "E51" "  -0.87%  "
